# Daily attendance processing - 2025-12-28 11:29:46
#
# Normalise the "Recorded By" (column G) lists so that the system-level
# recorder ("System" / "admin@admin.com") is always listed first, ahead of
# any human recorder email that got logged before it.
#
# Rule: for a comma-separated "Recorded By" value, if the 2nd entry is
# exactly "System" or "admin@admin.com" and the 1st entry is not already
# one of those, swap the first two entries. Everything else (3+ entries,
# single entries, already-correct ordering) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$colG = 7

$changedCount = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -ge 2) {
            $p0 = $parts[0]
            $p1 = $parts[1]

            $p0IsPriority = ($p0.Equals("System") -or $p0.Equals("admin@admin.com"))
            $p1IsPriority = ($p1.Equals("System") -or $p1.Equals("admin@admin.com"))

            if ((-not $p0IsPriority) -and $p1IsPriority) {
                $parts[0] = $p1
                $parts[1] = $p0
                $newval = $parts -join ", "
                $cell.Value = $newval
                $changedCount = $changedCount + 1
            }
        }
    }
}

Write-Host "Reordered Recorded By values in $changedCount row(s)."
